# Commit: "Researching About page inpage nav"
#
# Changes applied:
#   - D3 ("Time spent"): "50 minutes" -> "1hr 5min"
#   - E3 ("Notes"): the un-struck-through portion of the rich-text note is
#     expanded with an extra task, while the struck-through first run
#     ("Research mobile applications;") is left as-is.
#   - The active selection moves from B3 to D3.
#   - Column E is widened and no longer auto "best fit".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D3: time spent ---------------------------------------------------------
$ws.Range("D3").Value = "1hr 5min"

# --- E3: rich-text notes -----------------------------------------------------
$cell = $ws.Range("E3")

# Figure out where the struck-through lead-in run ends, using the text that
# is currently there, so we only touch the second run's content.
$oldFull = $cell.Text
$oldSecondRun = " determine PRL-appropriate structure; maybe ask Federica for opinion"
$splitAt = $oldFull.IndexOf($oldSecondRun)
$firstRun = $oldFull.Substring(0, $splitAt)

$newSecondRun = " Determine PRL-appropriate structure; Design XD mockup; Maybe ask Federica for opinion"

# Re-assign the full text, then reapply the strikethrough formatting to the
# first run only (the second run keeps its normal, non-struck-through look).
$cell.Value = $firstRun + $newSecondRun

$len1 = $firstRun.Length
$run1 = $cell.Characters(1, $len1)
$run1.Font.Strikethrough = $true

$len2 = $newSecondRun.Length
$run2 = $cell.Characters($len1 + 1, $len2)
$run2.Font.Strikethrough = $false

# --- Selection: B3 -> D3 ------------------------------------------------------
$null = $ws.Range("D3").Select()

# --- Column E width -----------------------------------------------------------
# Explicit width replaces the previous "best fit" auto-sizing.
$ws.Columns.Item(5).ColumnWidth = 36.3

Write-Host "Applied: D3/E3 text updates, selection D3, column E width."
